# Word COM-interop script implementing:
#   - lower-case the duplicated "Survey"/"Book" mentions (run-split preserved)
#   - move the "book author" run up so it ends its paragraph
#   - add a new "project" bullet (with its two sub-bullets) describing the
#     new "project" property that users can develop
#   - keep the old trailing (now-empty) ListParagraph paragraph, moving the
#     _GoBack bookmark into it
#
# Strategy: whole-paragraph replace via Range.InsertXML(fullParagraphXml)
# for paragraphs whose content changes, and Range.InsertParagraphAfter()
# to splice in brand-new paragraphs, which is then filled the same way.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Create management screen for Survey, a Survey has the following
#    information:" -> lower-case both "Survey" occurrences.
# ---------------------------------------------------------------------------
$pSurvey = $d.Paragraphs(3)
$xmlSurvey = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Create management screen for </w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">urvey, a </w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">urvey has the following </w:t></w:r>' +
  '<w:r><w:t>info</w:t></w:r>' +
  '<w:r><w:t>r</w:t></w:r>' +
  '<w:r><w:t>mation</w:t></w:r>' +
  '<w:r><w:t>:</w:t></w:r>' +
  '</w:p>'
$pSurvey.Range.InsertXML($xmlSurvey)

# ---------------------------------------------------------------------------
# 2) "Create management screen for Book, a Book has the following
#    information:" -> lower-case both "Book" occurrences.
# ---------------------------------------------------------------------------
$pBook = $d.Paragraphs(6)
$xmlBook = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Create management screen for </w:t></w:r>' +
  '<w:r><w:t>b</w:t></w:r>' +
  '<w:r><w:t>ook</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">, a </w:t></w:r>' +
  '<w:r><w:t>b</w:t></w:r>' +
  '<w:r><w:t>ook</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> has the following </w:t></w:r>' +
  '<w:r><w:t>information</w:t></w:r>' +
  '<w:r><w:t>:</w:t></w:r>' +
  '</w:p>'
$pBook.Range.InsertXML($xmlBook)

# ---------------------------------------------------------------------------
# 3) "An user is the book author" paragraph: drop the _GoBack bookmark here
#    (it moves to the trailing empty paragraph below) and let the
#    "book author" run simply end the paragraph.
# ---------------------------------------------------------------------------
$pAuthor = $d.Paragraphs(8)
$xmlAuthor = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t>A</w:t></w:r>' +
  '<w:r><w:t>n</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> user is </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">the </w:t></w:r>' +
  '<w:r><w:t>book author</w:t></w:r>' +
  '</w:p>'
$pAuthor.Range.InsertXML($xmlAuthor)

# ---------------------------------------------------------------------------
# 4) Insert three brand-new paragraphs right after it:
#      - "Create management screen for project, a project has the
#         following information:"
#      - " Name, Start Date, End Date, Customer Name, Status "
#      - "A list of users that will develop this project"
# ---------------------------------------------------------------------------
$pAuthor = $d.Paragraphs(8)
$pAuthor.Range.InsertParagraphAfter() | Out-Null
$pProject = $d.Paragraphs(9)
$xmlProject = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Create management screen for </w:t></w:r>' +
  '<w:r><w:t>project</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">, a </w:t></w:r>' +
  '<w:r><w:t>project</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> has the following information:</w:t></w:r>' +
  '</w:p>'
$pProject.Range.InsertXML($xmlProject)

$pProject = $d.Paragraphs(9)
$pProject.Range.InsertParagraphAfter() | Out-Null
$pFields = $d.Paragraphs(10)
$xmlFields = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve"> Name, Start Date, End Date,</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> Customer Name,</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> Status </w:t></w:r>' +
  '</w:p>'
$pFields.Range.InsertXML($xmlFields)

$pFields = $d.Paragraphs(10)
$pFields.Range.InsertParagraphAfter() | Out-Null
$pUsers = $d.Paragraphs(11)
$xmlUsers = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">A list of users that will </w:t></w:r>' +
  '<w:r><w:t>develop</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>this project</w:t></w:r>' +
  '</w:p>'
$pUsers.Range.InsertXML($xmlUsers)

# ---------------------------------------------------------------------------
# 5) The old trailing empty "ListParagraph" paragraph (no numPr) now picks
#    up the _GoBack bookmark that used to sit inside the author paragraph.
# ---------------------------------------------------------------------------
$pTail = $d.Paragraphs(12)
$xmlTail = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$pTail.Range.InsertXML($xmlTail)

Write-Output "done"
